$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header for column F
$ws.Range("F1").Value = "data"

# Update row 2 (João)
$ws.Range("C2").Value = 0.2
$ws.Range("E2").Value = 2000
$ws.Range("F2").Value = "2025-05-13 22:02:38"

# Update row 3 (Maria -> Claudia)
$ws.Range("A3").Value = "Claudia"
$ws.Range("B3").Value = 1900
$ws.Range("C3").Value = 0.5
$ws.Range("E3").Value = 1950
$ws.Range("F3").Value = "2025-05-13 22:02:50"

# Add new row 4 (Gabrielle)
$ws.Range("A4").Value = "Gabrielle"
$ws.Range("B4").Value = 1500
$ws.Range("C4").Value = 0.5
$ws.Range("D4").Value = 1000
$ws.Range("E4").Value = 1750
$ws.Range("F4").Value = "2025-05-13 22:03:00"
